$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("optimization_parameters")
$ws.Activate()

# Remove the redundant duplicated "value" header cells in C1:F1
$ws.Range("C1:F1").ClearContents()

# Insert a new row for the "L_curve" parameter right after the
# "Model"/production-function row (old row 8)
$ws.Rows.Item(9).Insert()

# Rename "Model" -> "production_function" label
$ws.Range("A8").Value = "production_function"

# Fill in the newly inserted row
$ws.Range("A9").Value = "L_curve"
$ws.Range("B9").Value = 1
$ws.Range("B9").NumberFormat = $ws.Range("B2").NumberFormat

# Remove the old "Deletion" parameter row (now shifted down to row 17
# after the insert above)
$ws.Rows.Item(17).Delete()

# The active cell on this sheet is now B10
$ws.Range("B10").Select()
